$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1366.6666
$ws.Range("I18").Value = 1050
$ws.Range("K18").Value = 1050
$ws.Range("M18").Value = -766

$ws.Range("H33").Value = 782.0909
$ws.Range("I33").Value = 355.8889
$ws.Range("K33").Value = 355.8889
$ws.Range("M33").Value = -126.8889

$ws.Range("H40").Value = 2579.4119
$ws.Range("I40").Value = 1980
$ws.Range("J40").Value = 2906.3635
$ws.Range("K40").Value = 1980
$ws.Range("L40").Value = 2906.3635
$ws.Range("M40").Value = -1805
$ws.Range("N40").Value = -3256.3635

$ws.Range("H125").Value = 360.8
$ws.Range("J125").Value = 422.5
$ws.Range("L125").Value = 3802.5
$ws.Range("N125").Value = -8722.5

$ws.Range("H132").Value = 3169.7385
$ws.Range("I132").Value = 3169.5574
$ws.Range("K132").Value = 9508.672200000001
$ws.Range("M132").Value = -6978.672200000001

$ws.Range("H137").Value = 3158.2058
$ws.Range("I137").Value = 2849.3333
$ws.Range("J137").Value = 3899.5
$ws.Range("K137").Value = 8547.999899999999
$ws.Range("L137").Value = 11698.5
$ws.Range("M137").Value = -5997.999899999999
$ws.Range("N137").Value = -16798.5

$ws.Range("H138").Value = 5391.018
$ws.Range("J138").Value = 5558.674
$ws.Range("L138").Value = 16676.022
$ws.Range("N138").Value = -26956.022

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 993.35297
$ws.Range("I97").Value = 777.6429000000001
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 777.6429000000001
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -281.6429000000001
$ws.Range("N97").Value = -2992

$ws.Range("H102").Value = 3545.8
$ws.Range("I102").Value = 3545.8
$ws.Range("K102").Value = 3545.8
$ws.Range("M102").Value = -1923.8

$ws.Range("H121").Value = 82000
$ws.Range("J121").Value = 82000
$ws.Range("L121").Value = 82000
$ws.Range("N121").Value = -85494

$ws.Range("H132").Value = 3554.7222
$ws.Range("I132").Value = 3373.1875
$ws.Range("K132").Value = 10119.5625
$ws.Range("M132").Value = -7589.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 43666.668
$ws.Range("I97").Value = 43666.668
$ws.Range("K97").Value = 43666.668
$ws.Range("M97").Value = -42675.668

$ws.Range("H99").Value = 4608.8237
$ws.Range("I99").Value = 4521.875
$ws.Range("K99").Value = 4521.875
$ws.Range("M99").Value = -3023.875

$ws.Range("H107").Value = 442.16666
$ws.Range("I107").Value = 388.4
$ws.Range("J107").Value = 711
$ws.Range("K107").Value = 388.4
$ws.Range("L107").Value = 711
$ws.Range("M107").Value = 1531.6
$ws.Range("N107").Value = -4551

$ws.Range("H134").Value = 1606.8889
$ws.Range("I134").Value = 1606.8889
$ws.Range("K134").Value = 4820.6667
$ws.Range("M134").Value = -2285.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 6069.6665
$ws.Range("J15").Value = 5104.5
$ws.Range("L15").Value = 5104.5
$ws.Range("N15").Value = -5444.5

$ws.Range("H107").Value = 2250.889
$ws.Range("J107").Value = 2449.625
$ws.Range("L107").Value = 2449.625
$ws.Range("N107").Value = -6289.625

$ws.Range("H134").Value = 1362.8695
$ws.Range("J134").Value = 3066.3333
$ws.Range("L134").Value = 9198.999899999999
$ws.Range("N134").Value = -14268.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 154.66667
$ws.Range("I33").Value = 42.333332
$ws.Range("J33").Value = 379.33334
$ws.Range("K33").Value = 253.999992
$ws.Range("L33").Value = 2276.00004
$ws.Range("M33").Value = 29.00000800000001
$ws.Range("N33").Value = -2842.00004

$ws.Range("H46").Value = 40288.23
$ws.Range("I46").Value = 1204.3182
$ws.Range("J46").Value = 255249.75
$ws.Range("K46").Value = 3612.9546
$ws.Range("L46").Value = 765749.25
$ws.Range("M46").Value = -3521.9546
$ws.Range("N46").Value = -765931.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2186.6667
$ws.Range("I102").Value = 2279.4443
$ws.Range("J102").Value = 1630
$ws.Range("K102").Value = 2279.4443
$ws.Range("L102").Value = 1630
$ws.Range("M102").Value = -657.4443000000001
$ws.Range("N102").Value = -4874

$ws.Range("H126").Value = 4740.5
$ws.Range("I126").Value = 4104.5
$ws.Range("K126").Value = 12313.5
$ws.Range("M126").Value = -9843.5

$ws.Range("H132").Value = 6275.569
$ws.Range("I132").Value = 6080.282
$ws.Range("J132").Value = 6910.25
$ws.Range("K132").Value = 18240.846
$ws.Range("L132").Value = 20730.75
$ws.Range("M132").Value = -15710.846
$ws.Range("N132").Value = -25790.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4189.273
$ws.Range("I7").Value = 3897
$ws.Range("J7").Value = 4700.75
$ws.Range("K7").Value = 3897
$ws.Range("L7").Value = 4700.75
$ws.Range("M7").Value = -3785
$ws.Range("N7").Value = -4924.75

$ws.Range("H22").Value = 1631.1666
$ws.Range("I22").Value = 947
$ws.Range("K22").Value = 947
$ws.Range("M22").Value = -652

$ws.Range("H27").Value = 1631.1666
$ws.Range("I27").Value = 947
$ws.Range("K27").Value = 947
$ws.Range("M27").Value = -840

$ws.Range("H68").Value = 3109.5
$ws.Range("I68").Value = 3610.5
$ws.Range("K68").Value = 3610.5
$ws.Range("M68").Value = -2861.5

$ws.Range("H71").Value = 3109.5
$ws.Range("I71").Value = 3610.5
$ws.Range("K71").Value = 18052.5
$ws.Range("M71").Value = -14308.5

$ws.Range("H126").Value = 4189.273
$ws.Range("I126").Value = 3897
$ws.Range("J126").Value = 4700.75
$ws.Range("K126").Value = 11691
$ws.Range("L126").Value = 14102.25
$ws.Range("M126").Value = -9221
$ws.Range("N126").Value = -19042.25

$ws.Range("H132").Value = 4064.6667
$ws.Range("I132").Value = 3198.25
$ws.Range("J132").Value = 5797.5
$ws.Range("K132").Value = 9594.75
$ws.Range("L132").Value = 17392.5
$ws.Range("M132").Value = -7064.75
$ws.Range("N132").Value = -22452.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2832.6667
$ws.Range("I96").Value = 2832.6667
$ws.Range("K96").Value = 2832.6667
$ws.Range("M96").Value = -1459.6667

$ws.Range("H132").Value = 999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 2997
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -8057

$ws.Range("H133").Value = 90000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H135").Value = 88482.81
$ws.Range("J135").Value = 88482.81
$ws.Range("L135").Value = 88482.81
$ws.Range("N135").Value = -98622.81

$ws.Range("H136").Value = 22873.107
$ws.Range("I136").Value = 26584.652
$ws.Range("K136").Value = 79753.95599999999
$ws.Range("M136").Value = -77203.95599999999
